$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for columns I and J, reusing the same style as the
# existing header cells (e.g. H1) by copying formats over.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 values for rows 2-46
$iValues = @(4,6,9,5,6,7,6,6,8,9,7,9,8,7,6,8,7,11,5,7,7,8,6,7,7,5,7,10,8,6,8,6,7,7,10,8,7,7,8,8,8,8,5,7,3)
# IF values for rows 2-46
$jValues = @(5,6,9,6,7,7,7,6,8,9,7,9,8,7,6,8,7,11,6,7,8,8,7,7,8,6,8,11,8,6,8,6,7,7,10,8,7,8,8,8,8,8,5,7,3)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
